$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing row (328) down through the new rows (329-337)
$ws.Range("A328:I328").Copy() | Out-Null
$ws.Range("A329:I337").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 329
$ws.Range("A329").Value = 327
$ws.Range("B329").Value = 42794
$ws.Range("C329").Value = 'triptech@trimet.org'
$ws.Range("D329").Value = 'Compliment'
$ws.Range("E329").Value = '"EXCELLENT!!!"'
$ws.Range("F329").ClearContents() | Out-Null
$ws.Range("G329").ClearContents() | Out-Null
$ws.Range("H329").Value = 'http://trimet.org/#planner/results/from=775+NE+5TH+ST%2C+GRESHAM+97030%3A%3A45.501008%2C-122.423025&to=Hollywood%2FNE+42nd+Ave+Transit+Center%2C+Portland%3A%3A45.53328%2C-122.620636&mode=RAIL%2CTRAM%2CSUBWAY%2CFUNICULAR%2CGONDOLA%2CWALK&m=pm&walk=1260&arr=A'
$ws.Range("I329").Value = 'Weekday'

# Row 330
$ws.Range("A330").Value = 328
$ws.Range("B330").Value = 42793
$ws.Range("C330").Value = 'triptech@trimet.org'
$ws.Range("D330").Value = 'Complaint'
$ws.Range("E330").Value = 'Unhappy with trip plan'
$ws.Range("F330").Value = 'Scheduling issue'
$ws.Range("G330").ClearContents() | Out-Null
$ws.Range("H330").Value = 'https://trimet.org/#planner/results/from=SW+Lombard+%26+Allen+-+Stop+ID+3455%3A%3A45.47699%2C-122.800202&to=7881+SW+CAPITOL+HWY%3A%3A45.467623%2C-122.714682&m=am&walk=1260&arr=D'
$ws.Range("I330").Value = 'Weekday'

# Row 331
$ws.Range("A331").Value = 329
$ws.Range("B331").Value = 42797
$ws.Range("C331").Value = 'triptech@trimet.org'
$ws.Range("D331").Value = 'Complaint'
$ws.Range("E331").Value = 'Unhappy with trip plan'
$ws.Range("F331").Value = 'Scheduling issue'
$ws.Range("G331").ClearContents() | Out-Null
$ws.Range("H331").Value = 'https://trimet.org/#/planner/results/itin_num=2&from=19725%20RIVER%20RD,%20Gladstone::45.376446,-122.60503&to=SE%20122nd%20%26%20Powell%20N,%20Portland%20(Stop%20ID%206655)::45.497345,-122.537506&Walk=1260&Arr=A'
$ws.Range("I331").Value = 'Weekday'

# Row 332
$ws.Range("A332").Value = 330
$ws.Range("B332").Value = 42799
$ws.Range("C332").Value = 'triptech@trimet.org'
$ws.Range("D332").Value = 'Complaint'
$ws.Range("E332").Value = 'Unhappy with trip plan'
$ws.Range("F332").Value = 'Maximum walk distance too low'
$ws.Range("G332").ClearContents() | Out-Null
$ws.Range("H332").Value = 'http://trimet.org/#planner/results/from=1500+SW+5TH+AVE%2C+PORTLAND%3A%3A45.513468%2C-122.680474&to=8470+SW+OLESON+RD%2C+Portland&m=am&walk=1260&arr=A'
$ws.Range("I332").Value = 'Saturday'

# Row 333
$ws.Range("A333").Value = 331
$ws.Range("B333").Value = 42799
$ws.Range("C333").Value = 'triptech@trimet.org'
$ws.Range("D333").Value = 'Complaint'
$ws.Range("E333").Value = 'Can''t plan trip outside district'
$ws.Range("F333").ClearContents() | Out-Null
$ws.Range("G333").ClearContents() | Out-Null
$ws.Range("H333").Value = 'http://trimet.org/#planner/results/from=SW+Bull+Mountain+Rd+%26+Oregon+Route+99W%2C+Tigard%3A%3A45.414692%2C-122.791534&to=Vancouver%2C+WA&m=am&walk=840&arr=D'
$ws.Range("I333").Value = 'Weekday'

# Row 334
$ws.Range("A334").Value = 332
$ws.Range("B334").Value = 42801
$ws.Range("C334").Value = 'triptech@trimet.org'
$ws.Range("D334").Value = 'Information request'
$ws.Range("E334").Value = 'Wanted to know why MAX platforms have 2 stop IDs'
$ws.Range("F334").ClearContents() | Out-Null
$ws.Range("G334").ClearContents() | Out-Null
$ws.Range("H334").Value = 'https://trimet.org/#planner/results/from=Clackamas+Town+Ctr+%26+Clackamas+Town+Center+TC%2C+Clackamas+County%3A%3A45.43555%2C-122.56904&to=NE+7th+Ave+MAX+Station+W%2C+Portland+(Stop+ID+8375)%3A%3A45.53015%2C-122.65828&mode=RAIL%2CTRAM%2CSUBWAY%2CFUNICULAR%2CGONDOLA%2CWALK&m=am&walk=1260&arr=A'
$ws.Range("I334").Value = 'Weekday'

# Row 335
$ws.Range("A335").Value = 333
$ws.Range("B335").Value = 42801
$ws.Range("C335").Value = 'triptech@trimet.org'
$ws.Range("D335").Value = 'Feature request'
$ws.Range("E335").Value = 'Would like to save trips'
$ws.Range("F335").ClearContents() | Out-Null
$ws.Range("G335").ClearContents() | Out-Null
$ws.Range("H335").Value = 'http://trimet.org/#planner/results/from=8885+SW+CANYON+RD%2C+Portland%3A%3A45.49791%2C-122.768684&to=1511+SW+PARK+AVE%2C+Portland%3A%3A45.514206%2C-122.68472&m=pm&walk=1260&arr=A'
$ws.Range("I335").Value = 'Weekday'

# Row 336
$ws.Range("A336").Value = 334
$ws.Range("B336").Value = 42800
$ws.Range("C336").Value = 'triptech@trimet.org'
$ws.Range("D336").Value = 'Complaint'
$ws.Range("E336").Value = 'Dangerous trip plan'
$ws.Range("F336").Value = 'Routed on footway that is too rough to use'
$ws.Range("G336").ClearContents() | Out-Null
$ws.Range("H336").Value = 'http://trimet.org/#/planner/results/itin_num=2&from=3030'
$ws.Range("I336").Value = 'Weekday'

# Row 337
$ws.Range("A337").Value = 335
$ws.Range("B337").Value = 42804
$ws.Range("C337").Value = 'triptech@trimet.org'
$ws.Range("D337").Value = 'Compliment'
$ws.Range("E337").Value = '"Thanks again"'
$ws.Range("F337").ClearContents() | Out-Null
$ws.Range("G337").ClearContents() | Out-Null
$ws.Range("H337").Value = 'https://trimet.org/#planner/results/from=3508+NE+BROADWAY%2C+Portland%3A%3A45.53481%2C-122.628&to=3911+SE+MILWAUKIE+AVE%2C+Portland%3A%3A45.494194%2C-122.65303'
$ws.Range("I337").Value = 'Weekday'

Write-Output "done"